# Update countries & provincias Spain
# - Swap Venezuela / Guinea-Bisau rows (new data causes Venezuela to overtake
#   Guinea-Bisau in the case-count ranking)
# - Swap Islas Turcas y Caicos / Groenlandia rows (same reason)
# - Refresh the case-count figures for several countries
# - Bump the "last updated" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Estados Unidos (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B4").Value = 1725117
$ws.Range("C4").Value = 18891
$ws.Range("D4").Value = 478225
$ws.Range("E4").Value = 1146352
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 735
$ws.Range("H4").Value = 100540

# Row 47 - Argentina
$ws.Range("B47").Value = 13228
$ws.Range("C47").Value = 600
$ws.Range("D47").Value = 4167
$ws.Range("E47").Value = 8571
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 23
$ws.Range("H47").Value = 490

# Row 104 used to be Guinea-Bisau; it now becomes Venezuela with refreshed data
$ws.Range("A104").Value = "Venezuela"
$ws.Range("B104").Value = 1211
$ws.Range("C104").Value = 34
$ws.Range("D104").Value = 302
$ws.Range("E104").Value = 898
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 11

# Row 105 used to be Venezuela; it now becomes Guinea-Bisau with refreshed data
$ws.Range("A105").Value = "Guinea-Bisau"
$ws.Range("B105").Value = 1178
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 42
$ws.Range("E105").Value = 1129
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 7

# Row 121 - Uruguay
$ws.Range("B121").Value = 789
$ws.Range("C121").Value = 2
$ws.Range("D121").Value = 638
$ws.Range("E121").Value = 129
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 22

# Row 123 - Principado de Andorra
$ws.Range("D123").Value = 676
$ws.Range("E123").Value = 36

# Row 207 used to be Groenlandia; it now becomes Islas Turcas y Caicos
$ws.Range("A207").Value = "Islas Turcas y Caicos"
$ws.Range("D207").Value = 10
$ws.Range("H207").Value = 1

# Row 208 used to be Islas Turcas y Caicos; it now becomes Groenlandia
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 0

# Timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 01:35"
